$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Commands")
$ws1 = $wb.Worksheets.Item("Sliders")

# The "EVENT CUSTOM SLIDERS" help workbook documents each available Command.
# Four obsolete Artisan commands (showCurve, showExtraCurve, showEvents,
# showBackgroundEvents) occupied rows 96-99 of the "Commands" sheet; they
# were dropped from the app, so the corresponding rows are removed here.
# Deleting these rows shifts everything below (RC Command, WebSocket
# Command, ...) up by four rows and shrinks the sheet from A1:C115 to
# A1:C111.
$ws2.Rows("96:99").Delete()

# Restore the cursor/selection state roughly as it was left by the editor:
# the "Sliders" sheet keeps its previous selection on B6 ...
$ws1.Range("B6").Select()

# ... while "Commands" remains the active tab, with the cursor now on the
# row that used to be just above the deleted block.
$ws2.Activate()
$ws2.Rows("95:95").Select()
